$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5414
$ws1.Range("F7").Value = 891
$ws1.Range("F9").Value = 2408
$ws1.Range("F12").Value = 2255
$ws1.Range("F13").Value = 68

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5414
$ws4.Range("F9").Value = 891
$ws4.Range("F11").Value = 2408
$ws4.Range("F15").Value = 2255
$ws4.Range("F16").Value = 68
